# template_import_parents_students.xlsx
# - Remove the "Email" (student's own email) column and its sample data
#   (HS001@email.com / HS002@email.com), including the hyperlink that was
#   attached to that sample value.
# - Insert a new "Năm học" (Academic Year) column right after "Lớp",
#   populated with sample academic-year ranges.
# - Re-center the "Khối" header.
# - Resize columns to their new best-fit widths.
# - Move the active selection to I7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the hyperlink that lived on B2 (mailto:HS001@email.com) before the
#    column that carries it disappears.
$ws.Range("B2").Hyperlinks.Delete()

# 2. Delete column B ("Email" header + HS001@email.com / HS002@email.com).
#    Everything to the right shifts one column left:
#    C->B (Ngay sinh), D->C (Gioi tinh), E->D (Khoi), F->E (Lop),
#    G->F (Ten phu huynh), H->G (SDT phu huynh), I->H (Email phu huynh)
$ws.Columns("B").Delete()

# 3. Insert a new blank column at F (right after the new "Lop" column E),
#    which will become "Nam hoc". This shifts the parent-info columns back
#    to the right: F->G (Ten phu huynh), G->H (SDT phu huynh), H->I (Email phu huynh)
$ws.Columns("F").Insert()

# 4. Populate the new "Nam hoc" column. Insert F3 before F2 so the shared
#    strings land in the same order as the authored workbook.
$ws.Range("F1").Value = "Năm học"
$ws.Range("F3").Value = "2020-2025"
$ws.Range("F2").Value = "2025-2030"

# 5. Re-center the "Khối" header (now D1).
$ws.Range("D1").HorizontalAlignment = -4108

# 6. Columns A-E, H, I, K already carry the right widths forward from the
#    delete/insert above. Only the brand new "Nam hoc" column (F) and the
#    "Ten phu huynh" column (G, squeezed by the new neighbour) need their
#    best-fit widths refreshed.
$ws.Columns("F").ColumnWidth = 8.877604166666666
$ws.Columns("G").ColumnWidth = 14.736979166666666

# 7. Move the selection like the authored session left it.
$ws.Range("I7").Select()
